$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("E2").Value = "No Responce"
$ws.Range("D2").Value = "2023-31-21 11:31:43"
$ws.Range("E3").Value = "No Responce"
$ws.Range("D3").Value = "2023-31-21 11:31:43"
$ws.Range("E4").Value = "No Responce"
$ws.Range("D4").Value = "2023-31-21 11:31:43"
